# Applies the "fixed some errors on Medium and NamingConvention" commit.
#
# Summary of the change:
#  - On the "@prefix" sheet, the "unitLength" and "image" prefix rows swap
#    places (unitLength now row 15, image now row 16) and three brand new
#    prefix rows are appended: pixels, channel, bindata - each pointing at
#    a new namespace URI under .../samples/multi-channel/<name>/
#  - Everywhere the old bracketed placeholder-style local identifiers
#    ([pixels:image0], [channel:0..2], [bindata:0..2]) were used as values,
#    they are replaced by real prefixed identifiers that use the new
#    prefixes (pixels:pixels0, channel:channel0..2, bindata:bindata0..2).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "@prefix" sheet: reorder unitLength/image rows and add 3 new prefixes
# ---------------------------------------------------------------------
$wsPrefix = $wb.Worksheets.Item("@prefix")

$wsPrefix.Cells.Item(15, 1).Value = "unitLength"
$wsPrefix.Cells.Item(15, 2).Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/Unit/UnitLength#"

$wsPrefix.Cells.Item(16, 1).Value = "image"
$wsPrefix.Cells.Item(16, 2).Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/multi-channel/image/"

$wsPrefix.Cells.Item(17, 1).Value = "pixels"
$wsPrefix.Cells.Item(17, 2).Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/multi-channel/pixels/"

$wsPrefix.Cells.Item(18, 1).Value = "channel"
$wsPrefix.Cells.Item(18, 2).Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/multi-channel/channel/"

$wsPrefix.Cells.Item(19, 1).Value = "bindata"
$wsPrefix.Cells.Item(19, 2).Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/multi-channel/bindata/"

# ---------------------------------------------------------------------
# 2. "Image" sheet: local identifier now uses the "pixels" prefix
# ---------------------------------------------------------------------
$wsImage = $wb.Worksheets.Item("Image")
$wsImage.Range("E5").Value = "pixels:pixels0"

# ---------------------------------------------------------------------
# 3. "Pixels" sheet: same pixels identifier, plus channel/bindata values
# ---------------------------------------------------------------------
$wsPixels = $wb.Worksheets.Item("Pixels")
$wsPixels.Range("B5").Value = "pixels:pixels0"
$wsPixels.Range("B6").Value = "pixels:pixels0"
$wsPixels.Range("B7").Value = "pixels:pixels0"

$wsPixels.Range("M5").Value = "channel:channel0"
$wsPixels.Range("N5").Value = "bindata:bindata0"
$wsPixels.Range("M6").Value = "channel:channel1"
$wsPixels.Range("N6").Value = "bindata:bindata1"
$wsPixels.Range("M7").Value = "channel:channel2"
$wsPixels.Range("N7").Value = "bindata:bindata2"

$wsPixels.Columns.Item(2).ColumnWidth = 12.0
$wsPixels.Columns.Item(13).ColumnWidth = 15.428571428571429
$wsPixels.Columns.Item(14).ColumnWidth = 14.714285714285715

# ---------------------------------------------------------------------
# 4. "Channel" sheet: local identifiers now use the "channel" prefix
# ---------------------------------------------------------------------
$wsChannel = $wb.Worksheets.Item("Channel")
$wsChannel.Range("B5").Value = "channel:channel0"
$wsChannel.Range("B6").Value = "channel:channel1"
$wsChannel.Range("B7").Value = "channel:channel2"

$wsChannel.Columns.Item(2).ColumnWidth = 15.428571428571429

# ---------------------------------------------------------------------
# 5. "Binary_Data" sheet: local identifiers now use the "bindata" prefix
# ---------------------------------------------------------------------
$wsBinData = $wb.Worksheets.Item("Binary_Data")
$wsBinData.Range("B5").Value = "bindata:bindata0"
$wsBinData.Range("B6").Value = "bindata:bindata1"
$wsBinData.Range("B7").Value = "bindata:bindata2"

$wsBinData.Columns.Item(2).ColumnWidth = 14.714285714285715
